$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.342.27"
$ws.Range("E2").Value = "  -2.74%  "
$ws.Range("D3").Value = "1.940.97"
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7224"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.21%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -4.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "29.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07385"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8161"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08121"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").Value = "1.937.18"
$ws.Range("E13").Value = "  -2.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.490"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("D17").Value = "30.351.34"
$ws.Range("E17").Value = "  -2.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008331"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.886"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "2.193.46"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.856"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.406"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("E28").Value = "  -2.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1322"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.32%  "
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.346"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.463"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.247"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05268"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.308"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7561"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01995"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.857"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "81.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.620"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4562"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.044"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8471"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.862"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.518"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4195"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.505"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.89%  "
